# iecodebook: update run files
# The "survey" sheet lists, for each dataset variable, the current
# name/label (columns E/F) and an optional new name/label (columns A/B).
# Here we add a new entry: keep the "price" variable's name but relabel
# it to "Cost".

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("survey")
$ws.Activate()

$ws.Range("A4").Value = "price"
$ws.Range("B4").Value = "Cost"

$ws.Range("B5").Select()
